# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Sargatanas_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 512.3
$ws.Range("I28").Value = 553.6667
$ws.Range("K28").Value = 553.6667
$ws.Range("M28").Value = -68.66669999999999
# Row 58
$ws.Range("H58").Value = 841.36365
$ws.Range("I58").Value = 234.57143
$ws.Range("J58").Value = 1903.25
$ws.Range("K58").Value = 703.71429
$ws.Range("L58").Value = 5709.75
$ws.Range("M58").Value = -553.71429
$ws.Range("N58").Value = -6009.75
# Row 61
$ws.Range("H61").Value = 1350
$ws.Range("I61").Value = 1300
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 3900
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -3728
$ws.Range("N61").Value = -4544
# Row 98
$ws.Range("H98").Value = 4212.174
$ws.Range("I98").Value = 4212.174
$ws.Range("K98").Value = 4212.174
$ws.Range("M98").Value = -2714.174
# Row 99
$ws.Range("H99").Value = 222.875
$ws.Range("I99").Value = 227
$ws.Range("K99").Value = 681
$ws.Range("M99").Value = 817
# Row 100
$ws.Range("H100").Value = 3442.4
$ws.Range("I100").Value = 3419.1428
$ws.Range("J100").Value = 3496.6667
$ws.Range("K100").Value = 3419.1428
$ws.Range("L100").Value = 3496.6667
$ws.Range("M100").Value = -2878.1428
$ws.Range("N100").Value = -4578.6667
# Row 101
$ws.Range("H101").Value = 782.3333
$ws.Range("I101").Value = 594.2857
$ws.Range("K101").Value = 1782.8571
$ws.Range("M101").Value = -160.8571000000002
# Row 104
$ws.Range("H104").Value = 951.4
$ws.Range("I104").Value = 951.4
$ws.Range("K104").Value = 2854.2
$ws.Range("M104").Value = -1107.2
# Row 106
$ws.Range("H106").Value = 3316.0908
$ws.Range("I106").Value = 3347.8
$ws.Range("K106").Value = 3347.8
$ws.Range("M106").Value = -2716.8
# Row 111
$ws.Range("H111").Value = 7355726.5
$ws.Range("I111").Value = 15626264
$ws.Range("J111").Value = 4137.778
$ws.Range("K111").Value = 46878792
$ws.Range("L111").Value = 12413.334
$ws.Range("M111").Value = -46875725
$ws.Range("N111").Value = -18547.334
# Row 113
$ws.Range("H113").Value = 21751608
$ws.Range("I113").Value = 2642.5557
$ws.Range("J113").Value = 35733090
$ws.Range("K113").Value = 2642.5557
$ws.Range("L113").Value = 35733090
$ws.Range("M113").Value = 611.4443000000001
$ws.Range("N113").Value = -35739598
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 115
$ws.Range("H115").Value = 1816.6666
$ws.Range("I115").Value = 1500
$ws.Range("J115").Value = 1975
$ws.Range("K115").Value = 4500
$ws.Range("L115").Value = 5925
$ws.Range("M115").Value = -2933
$ws.Range("N115").Value = -9059
# Row 118
$ws.Range("H118").Value = 1703.8
$ws.Range("I118").Value = 1703.8
$ws.Range("K118").Value = 5111.4
$ws.Range("M118").Value = -3454.4
# Row 122
$ws.Range("H122").Value = 4212.174
$ws.Range("I122").Value = 4212.174
$ws.Range("K122").Value = 12636.522
$ws.Range("M122").Value = -10186.522
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 127
$ws.Range("H127").Value = 759
$ws.Range("I127").Value = 759
$ws.Range("K127").Value = 2277
$ws.Range("M127").Value = 2683
# Row 129
$ws.Range("H129").Value = 1913.6154
$ws.Range("I129").Value = 1459.3334
$ws.Range("K129").Value = 4378.0002
$ws.Range("M129").Value = 621.9997999999996
# Row 132
$ws.Range("H132").Value = 1064.5312
$ws.Range("I132").Value = 1084.3334
$ws.Range("K132").Value = 3253.0002
$ws.Range("M132").Value = -723.0001999999999
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 138
$ws.Range("H138").Value = 2515.62
$ws.Range("I138").Value = 1997.0741
$ws.Range("J138").Value = 2707.411
$ws.Range("K138").Value = 5991.2223
$ws.Range("L138").Value = 8122.233
$ws.Range("M138").Value = -851.2223000000004
$ws.Range("N138").Value = -18402.233

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2006278.9
$ws.Range("I32").Value = 2227498.8
$ws.Range("K32").Value = 2227498.8
$ws.Range("M32").Value = -2227211.8
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 122
$ws.Range("H122").Value = 3414.4375
$ws.Range("I122").Value = 2389.4546
$ws.Range("K122").Value = 7168.3638
$ws.Range("M122").Value = -4718.3638
# Row 132
$ws.Range("H132").Value = 6081.778
$ws.Range("I132").Value = 3080.1428
$ws.Range("J132").Value = 7991.909
$ws.Range("K132").Value = 9240.428400000001
$ws.Range("L132").Value = 23975.727
$ws.Range("M132").Value = -6710.428400000001
$ws.Range("N132").Value = -29035.727
# Row 135
$ws.Range("H135").Value = 95143
$ws.Range("J135").Value = 95143
$ws.Range("L135").Value = 95143
$ws.Range("N135").Value = -105283

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3421
$ws.Range("I105").Value = 2699.8333
$ws.Range("K105").Value = 2699.8333
$ws.Range("M105").Value = -952.8332999999998
# Row 107
$ws.Range("H107").Value = 86540440
$ws.Range("I107").Value = 86540440
$ws.Range("K107").Value = 86540440
$ws.Range("M107").Value = -86538520
# Row 134
$ws.Range("H134").Value = 6103336.5
$ws.Range("I134").Value = 9618079
$ws.Range("K134").Value = 28854237
$ws.Range("M134").Value = -28851702

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4346.5454
$ws.Range("I16").Value = 3779.9
$ws.Range("K16").Value = 3779.9
$ws.Range("M16").Value = -3492.9
# Row 31
$ws.Range("H31").Value = 8939.617
$ws.Range("I31").Value = 2674.1667
$ws.Range("J31").Value = 12357.137
$ws.Range("K31").Value = 2674.1667
$ws.Range("L31").Value = 12357.137
$ws.Range("M31").Value = -2379.1667
$ws.Range("N31").Value = -12947.137
# Row 34
$ws.Range("H34").Value = 8939.617
$ws.Range("I34").Value = 2674.1667
$ws.Range("J34").Value = 12357.137
$ws.Range("K34").Value = 2674.1667
$ws.Range("L34").Value = 12357.137
$ws.Range("M34").Value = -2472.1667
$ws.Range("N34").Value = -12761.137
# Row 113
$ws.Range("H113").Value = 4346.5454
$ws.Range("I113").Value = 3779.9
$ws.Range("K113").Value = 3779.9
$ws.Range("M113").Value = -1609.9

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 760.6667
$ws.Range("I63").Value = 760.6667
$ws.Range("K63").Value = 2282.0001
$ws.Range("M63").Value = -1533.0001
# Row 66
$ws.Range("H66").Value = 760.6667
$ws.Range("I66").Value = 760.6667
$ws.Range("K66").Value = 6846.0003
$ws.Range("M66").Value = -3102.0003
# Row 140
$ws.Range("H140").Value = 77151.63
$ws.Range("I140").Value = 126716.375
$ws.Range("K140").Value = 380149.125
$ws.Range("M140").Value = -374969.125

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 10242.571
$ws.Range("J35").Value = 9000
$ws.Range("L35").Value = 9000
$ws.Range("N35").Value = -9596
# Row 93
$ws.Range("H93").Value = 69990
$ws.Range("J93").Value = 69990
$ws.Range("L93").Value = 69990
$ws.Range("N93").Value = -73734
# Row 113
$ws.Range("H113").Value = 5280.769
$ws.Range("I113").Value = 2948.1667
$ws.Range("J113").Value = 7280.143
$ws.Range("K113").Value = 2948.1667
$ws.Range("L113").Value = 7280.143
$ws.Range("M113").Value = -778.1667000000002
$ws.Range("N113").Value = -11620.143
# Row 122
$ws.Range("H122").Value = 10351472
$ws.Range("I122").Value = 24144188
$ws.Range("K122").Value = 72432564
$ws.Range("M122").Value = -72430114
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 132
$ws.Range("H132").Value = 7431.4614
$ws.Range("I132").Value = 2802
$ws.Range("J132").Value = 12832.5
$ws.Range("K132").Value = 8406
$ws.Range("L132").Value = 38497.5
$ws.Range("M132").Value = -5876
$ws.Range("N132").Value = -43557.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 878.2778
$ws.Range("J16").Value = 70
$ws.Range("L16").Value = 70
$ws.Range("N16").Value = -410
# Row 61
$ws.Range("H61").Value = 4196.381
$ws.Range("I61").Value = 1755.8462
$ws.Range("J61").Value = 8162.25
$ws.Range("K61").Value = 1755.8462
$ws.Range("L61").Value = 8162.25
$ws.Range("M61").Value = -1553.8462
$ws.Range("N61").Value = -8566.25
# Row 68
$ws.Range("H68").Value = 5417
$ws.Range("I68").Value = 3435.5
$ws.Range("J68").Value = 7681.5713
$ws.Range("K68").Value = 3435.5
$ws.Range("L68").Value = 7681.5713
$ws.Range("M68").Value = -2686.5
$ws.Range("N68").Value = -9179.5713
# Row 71
$ws.Range("H71").Value = 5417
$ws.Range("I71").Value = 3435.5
$ws.Range("J71").Value = 7681.5713
$ws.Range("K71").Value = 17177.5
$ws.Range("L71").Value = 38407.85649999999
$ws.Range("M71").Value = -13433.5
$ws.Range("N71").Value = -45895.85649999999
# Row 82
$ws.Range("H82").Value = 2631.2222
$ws.Range("I82").Value = 1051.3334
$ws.Range("J82").Value = 3421.1667
$ws.Range("K82").Value = 1051.3334
$ws.Range("L82").Value = 3421.1667
$ws.Range("M82").Value = -690.3334
$ws.Range("N82").Value = -4143.1667
# Row 85
$ws.Range("H85").Value = 2631.2222
$ws.Range("I85").Value = 1051.3334
$ws.Range("J85").Value = 3421.1667
$ws.Range("K85").Value = 1051.3334
$ws.Range("L85").Value = 3421.1667
$ws.Range("M85").Value = 196.6666
$ws.Range("N85").Value = -5917.1667
# Row 113
$ws.Range("H113").Value = 4196.381
$ws.Range("I113").Value = 1755.8462
$ws.Range("J113").Value = 8162.25
$ws.Range("K113").Value = 1755.8462
$ws.Range("L113").Value = 8162.25
$ws.Range("M113").Value = 414.1538
$ws.Range("N113").Value = -12502.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 380318.9
$ws.Range("I122").Value = 1369602.6
$ws.Range("J122").Value = 9337.5
$ws.Range("K122").Value = 4108807.8
$ws.Range("L122").Value = 28012.5
$ws.Range("M122").Value = -4106357.8
$ws.Range("N122").Value = -32912.5

Write-Host "Applied updates: 252 cell writes, 5 cell clears"
